$p = $ppt.ActivePresentation

foreach ($idx in 14, 15, 16) {
    $s = $p.Slides.Item($idx)
    $shape = $s.Shapes.Item(1)
    $tbl = $shape.Table
    $tbl.ApplyStyle("{871F94BB-F7C7-445F-992E-DA2D23663825}")
}
